# Auto-generated edit script applying the cryptos.xlsx data refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.725.95"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "3.800.24"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'705.43"
$ws.Range("E5").Value = "  +1.34%  "
$ws.Range("D6").Value = "'170.05"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("D7").Value = "3.802.13"
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").Value = "'0.160"
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("D11").Value = "'7.36"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").Value = "'0.456"
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("E13").Value = "  -2.22%  "
$ws.Range("D14").Value = "'36.04"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").Value = "4.439.96"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("D16").Value = "3.799.13"
$ws.Range("E16").Value = "  -1.49%  "
$ws.Range("D17").Value = "70.733.76"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("E19").Value = "  -1.68%  "
$ws.Range("D20").Value = "'17.40"
$ws.Range("E20").Value = "  -2.05%  "
$ws.Range("D21").Value = "'497.46"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("D22").Value = "'10.60"
$ws.Range("E22").Value = "  -5.06%  "
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("D24").Value = "'84.58"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("E25").Value = "  -1.32%  "
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("D27").Value = "'10.42"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D28").Value = "3.949.86"
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  -4.79%  "
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("E32").Value = "  -4.20%  "
$ws.Range("E33").Value = "  -3.89%  "
$ws.Range("D34").Value = "'29.01"
$ws.Range("E34").Value = "  -2.40%  "
$ws.Range("E35").Value = "  -3.44%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'9.10"
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("B38").Value = "RenzoRestakedETH"
$ws.Range("C38").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D38").Value = "3.769.48"
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("E39").Value = "  -3.57%  "
$ws.Range("E40").Value = "  -1.82%  "
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("E43").Value = "  -4.32%  "
$ws.Range("D44").Value = "'0.000327"
$ws.Range("E44").Value = "  +6.73%  "
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").Value = "'164.87"
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("D48").Value = "'426.36"
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("E51").Value = "  -1.60%  "
